# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" sheet right after "总计", shifting the
#    existing quarter sheets (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3) later.
# 2) Update the "总计" (totals) sheet with a new top data row for 2022-Q3
#    and renumber the existing index column accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: update the "总计" summary sheet
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Give the new last row (A6) the same bold/centered/bordered look as the
# rest of the index column by copying the formatting from the cell above,
# before any values are written.
$totals.Range("A5").Copy() | Out-Null
$totals.Range("A6").PasteSpecial(-4122) | Out-Null

# Rewrite the full data block (rows 2-6) top to bottom with the new
# 2022-Q3 row inserted first and everything else pushed down one slot.
# Values are written as literals (rather than copied cell-to-cell) to
# avoid COM variant round-tripping from perturbing the stored doubles.
$rows = @(
    @("2022-Q3", 2, 2.22),
    @("2022-Q2", 7, 3.87),
    @("2022-Q1", 15, 6.37),
    @("2021-Q4", 5, 1.4),
    @("2021-Q3", 1, 0.02)
)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $totals.Cells.Item($r, 1).Value = $i
    $totals.Cells.Item($r, 2).Value = $data[0]
    $totals.Cells.Item($r, 3).Value = $data[1]
    $totals.Cells.Item($r, 4).Value = $data[2]
}

# ---------------------------------------------------------------------
# Part 2: create the new "2022-Q3" detail sheet right after "总计"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($afterSheet.Next)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data rows: fund code/name/size/position/ratio/value are stored as text
# (to preserve formatting such as leading zeros and trailing decimals),
# only the rank column is numeric.
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "'519702"
$q3.Cells.Item(2, 3).Value = "'交银趋势优先混合A"
$q3.Cells.Item(2, 4).Value = "'83.94"
$q3.Cells.Item(2, 5).Value = "'81.61"
$q3.Cells.Item(2, 6).Value = "'2.30"
$q3.Cells.Item(2, 7).Value = "'1.9306"
$q3.Cells.Item(2, 8).Value = 8

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "'013430"
$q3.Cells.Item(3, 3).Value = "'交银趋势优先混合C"
$q3.Cells.Item(3, 4).Value = "'12.61"
$q3.Cells.Item(3, 5).Value = "'81.61"
$q3.Cells.Item(3, 6).Value = "'2.30"
$q3.Cells.Item(3, 7).Value = "'0.2900"
$q3.Cells.Item(3, 8).Value = 8

# Formatting: bold, centered, thin-bordered header row and index column,
# matching the look of the other quarter sheets.
$hdr = $q3.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

$idx = $q3.Range("A2:A3")
$idx.Font.Bold = $true
$idx.Borders.LineStyle = 1
$idx.HorizontalAlignment = -4108
$idx.VerticalAlignment = -4160

Write-Output "2022-Q3 sheet added and totals updated"
